# Remove the trailing "Ver no Jupiter ..." / "© 2020 ..." footer block
# (plus the blank paragraph that precedes it) that the Jekyll site build
# appended after the bibliography's last reference.

$d = $word.ActiveDocument

# Locate the "Ver no Jupiter" paragraph.
$rVer = $d.Content
$null = $rVer.Find.Execute("Ver no Jupiter Salvar em pdf Salvar em docx", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$verParaIndex = $rVer.Paragraphs.Item(1).Index

# Locate the "© 2020 ... Powered by Jekyll ..." paragraph.
$rCopy = $d.Content
$null = $rCopy.Find.Execute("Powered by Jekyll and Github pages", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$copyParaIndex = $rCopy.Paragraphs.Item(1).Index

# The blank paragraph right before "Ver no Jupiter" is also removed.
$blankPara = $d.Paragraphs.Item($verParaIndex - 1)
$copyPara = $d.Paragraphs.Item($copyParaIndex)

$deleteRange = $d.Range($blankPara.Range.Start, $copyPara.Range.End)
$deleteRange.Delete()
